$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# YDS sheet: append newly-logged (Week 15) / simulated (Week 16) play
# yardage values to the running per-play logs.
# ---------------------------------------------------------------------
$ydsWs = $wb.Worksheets.Item("YDS")

$ydsWs.Range("B2").Value2 = $ydsWs.Range("B2").Value2 + " 2 1 10 3 0 7 2 7 0 0 1 2 3 -3 6 -5"
$ydsWs.Range("C2").Value2 = $ydsWs.Range("C2").Value2 + " 3 6 2 1 3 8 9 2 1 2 3 6 0 3 9 1 2 0 9 12 1 6 3 10 7 2 4 4 15 4 5 2 11 0 20 6 7 3 1 4 3"
$ydsWs.Range("B3").Value2 = $ydsWs.Range("B3").Value2 + " 16 4 6 3 -6 11 14 7 19 5 5 17 12 7 5 8 4"
$ydsWs.Range("C3").Value2 = $ydsWs.Range("C3").Value2 + " 7 5 11 12 1 4 8 5 18 7 9 13 11 5 -4 -3 3 5 12 5 9 8"

# ---------------------------------------------------------------------
# OFF sheet: season-to-date offensive totals (Home row 2, Road row 3)
# ---------------------------------------------------------------------
$offWs = $wb.Worksheets.Item("OFF")

$offWs.Range("B2").Value = 6
$offWs.Range("C2").Value = 159
$offWs.Range("D2").Value = 11
$offWs.Range("G2").Value = 40
$offWs.Range("J2").Value = 19
$offWs.Range("L2").Value = 302
$offWs.Range("M2").Value = 197
$offWs.Range("Q2").Value = 518

$offWs.Range("C3").Value = 161
$offWs.Range("E3").Value = 34
$offWs.Range("F3").Value = 106
$offWs.Range("G3").Value = 38
$offWs.Range("H3").Value = 31
$offWs.Range("I3").Value = 59
$offWs.Range("J3").Value = 67
$offWs.Range("N3").Value = 16

# ---------------------------------------------------------------------
# DEF sheet: season-to-date defensive totals (Home row 2, Road row 3)
# ---------------------------------------------------------------------
$defWs = $wb.Worksheets.Item("DEF")

$defWs.Range("B2").Value = 9
$defWs.Range("C2").Value = 186
$defWs.Range("D2").Value = 12
$defWs.Range("E2").Value = 15
$defWs.Range("F2").Value = 58
$defWs.Range("G2").Value = 48
$defWs.Range("H2").Value = 4
$defWs.Range("J2").Value = 29
$defWs.Range("L2").Value = 247
$defWs.Range("M2").Value = 165
$defWs.Range("O2").Value = 18
$defWs.Range("P2").Value = 10
$defWs.Range("Q2").Value = 510

$defWs.Range("C3").Value = 143
$defWs.Range("D3").Value = 3
$defWs.Range("E3").Value = 29
$defWs.Range("F3").Value = 82
$defWs.Range("G3").Value = 31
$defWs.Range("H3").Value = 38
$defWs.Range("I3").Value = 49
$defWs.Range("J3").Value = 38
$defWs.Range("N3").Value = 28

# ---------------------------------------------------------------------
# ST sheet: special-teams totals + per-kick logs
# ---------------------------------------------------------------------
$stWs = $wb.Worksheets.Item("ST")

$stWs.Range("B2").Value = 73
$stWs.Range("D2").Value = 62
$stWs.Range("F2").Value = 241
$stWs.Range("G2").Value = 227
$stWs.Range("J2").Value = 104
$stWs.Range("K2").Value = 99
$stWs.Range("L2").Value = 73
$stWs.Range("M2").Value = 59
$stWs.Range("N2").Value = 20

$stWs.Range("B3").Value = 37

$stWs.Range("B4").Value2 = $stWs.Range("B4").Value2 + " 64 59"
$stWs.Range("B5").Value2 = $stWs.Range("B5").Value2 + " 21 19"
$stWs.Range("B6").Value2 = $stWs.Range("B6").Value2 + " 20 18"
$stWs.Range("D3").Value2 = $stWs.Range("D3").Value2 + " 27 51 32 40 49"
$stWs.Range("D4").Value2 = $stWs.Range("D4").Value2 + " 0 55 0 0 0"
$stWs.Range("D5").Value2 = $stWs.Range("D5").Value2 + " 17 11 11"

# ---------------------------------------------------------------------
# TURNS sheet: turnover totals
# ---------------------------------------------------------------------
$turnsWs = $wb.Worksheets.Item("TURNS")

$turnsWs.Range("B2").Value = 6
$turnsWs.Range("C2").Value = 6
$turnsWs.Range("D2").Value = 6
$turnsWs.Range("E2").Value = 10

$turnsWs.Range("D3").Value = 5

# ---------------------------------------------------------------------
# PEN sheet: penalty totals
# ---------------------------------------------------------------------
$penWs = $wb.Worksheets.Item("PEN")

$penWs.Range("B2").Value = 18
$penWs.Range("D2").Value = 9
